# BOM.xlsx edit: fix wrong MPN / footprint entries and add spacing rows
# for production-data bookkeeping (per commit message).
#
# Net effect on the "BOM" sheet:
#   * Blank spacer rows are introduced between component groups
#     (Diodes -> Jacks -> Transistors -> Resistors), and within the
#     transistor group, so the existing rows 9-23 end up re-numbered to
#     11,14,15,22,24,28-37 while keeping their original content.
#   * The TO-220/TO-229 transistor footprint text is corrected from
#     "TO229P800X325X1420-3 (Bzw. TO220 in Kicad)" to
#     "TO229P780X300X1364-3" (used by Q1 and Q2).
#   * Minor cosmetic view changes (zoom level, selected cell, column A
#     width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert blank rows so the existing data shifts down into the new
#    row numbers, leaving gaps between the component groups.
#    (Row numbers below refer to positions *as they are at the moment
#    of each insert*, i.e. before any later inserts further down.)
# ---------------------------------------------------------------------

# Gap between the diode group (row 8) and the jack group (old row 9).
$ws.Range("A9:A10").EntireRow.Insert()

# Gap between the jack group and the transistor group (old row 10, Q1).
$ws.Range("A12:A13").EntireRow.Insert()

# Gap inside the transistor group, between Q12/Q13 and Q2.
$ws.Range("A16:A21").EntireRow.Insert()

# Gap between Q2 and Q3/Q9.
$ws.Range("A23:A23").EntireRow.Insert()

# Gap between the transistor group and the resistor group.
$ws.Range("A25:A27").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2) The inserted rows are completely empty, so clear them fully
#    (no leftover formatting/row element should remain for them).
# ---------------------------------------------------------------------
$blankRows = @(9, 10, 12, 13, 16, 17, 18, 19, 20, 21, 23, 25, 26, 27)
foreach ($r in $blankRows) {
    $ws.Rows.Item($r).Clear()
}

# ---------------------------------------------------------------------
# 3) Correct the erroneous MPN / footprint text for the TO-220 package
#    used by Q1 (now row 14) and Q2 (now row 22).
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 2).Value = "TO229P780X300X1364-3"
$ws.Cells.Item(22, 2).Value = "TO229P780X300X1364-3"

# ---------------------------------------------------------------------
# 4) Cosmetic view updates: narrower zoom, column A widened, and the
#    active selection moved to B9 (first of the newly shifted rows).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.25

$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("B9").Select()
